$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.548.21'
$ws.Range("E2").Value = '  +0.64%  '
$ws.Range("D3").Value = '3.252.05'
$ws.Range("E3").Value = '  +2.73%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '607.19'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.79%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '157.21'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +2.15%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").Value = '3.250.11'
$ws.Range("E8").Value = '  +2.70%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.549'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.36%  '
$ws.Range("E10").Value = '  +2.59%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.85'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +5.18%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.500'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -3.49%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000270'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +1.40%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '39.08'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +2.01%  '
$ws.Range("D15").Value = '3.785.01'
$ws.Range("E15").Value = '  +2.74%  '
$ws.Range("D16").Value = '66.612.56'
$ws.Range("E16").Value = '  +0.71%  '
$ws.Range("E17").Value = '  +0.60%  '
$ws.Range("D18").Value = '3.259.14'
$ws.Range("E18").Value = '  +3.02%  '
$ws.Range("E19").Value = '  +1.18%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '506.89'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.51%  '
$ws.Range("E21").Value = '  +0.37%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.751'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +3.21%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.14'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +1.36%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '14.72'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.37%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '87.07'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +3.02%  '
$ws.Range("E26").Value = '  -0.05%  '
$ws.Range("E27").Value = '  +1.33%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.15'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.89%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.40'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +0.71%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.131'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +48.52%  '
$ws.Range("E31").Value = '  -5.00%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.94'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -1.02%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '28.04'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +0.30%  '
$ws.Range("E34").Value = '  +0.08%  '
$ws.Range("E35").Value = '  -3.15%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.46'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.32%  '
$ws.Range("E37").Value = '  +20.73%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '55.66'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +1.80%  '
$ws.Range("D39").Value = '0.0₃0783'
$ws.Range("E39").Value = '  +16.01%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '494.20'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.86%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0422'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +0.81%  '
$ws.Range("E42").Value = '  -0.30%  '
$ws.Range("E43").Value = '  +1.02%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.293'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.64%  '
$ws.Range("E45").Value = '  +4.11%  '
$ws.Range("D46").Value = '2.986.96'
$ws.Range("E46").Value = '  +5.93%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '28.87'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +3.90%  '
$ws.Range("E48").Value = '  +5.87%  '
$ws.Range("E49").Value = '  +2.36%  '
$ws.Range("E50").Value = '  -0.04%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '121.13'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.47%  '
